$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.356.63"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "2.641.44"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.39%  "

$ws.Range("E10").Value = "  +0.80%  "

$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("D13").Value = "3.108.11"
$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("D14").Value = "59.355.30"
$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.35%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.694.08"
$ws.Range("E16").Value = "  +1.15%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "342.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.24%  "

$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("D26").Value = "2.766.13"

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.86%  "

$ws.Range("D29").Value = "0.0₃0795"
$ws.Range("E29").Value = "  -1.02%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  -4.38%  "

$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.42%  "

$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("E36").Value = "  -2.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.858"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.67%  "

$ws.Range("E40").Value = "  -2.04%  "

$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.600"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "269.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.20%  "

$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("D49").Value = "2.035.30"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.79%  "
